$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 751
$ws.Range("I2").Value = 1805
$ws.Range("J2").Value = 7614
$ws.Range("K2").Value = 54
$ws.Range("L2").Value = 2165
$ws.Range("M2").Value = 136
$ws.Range("N2").Value = 1362
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 30
$ws.Range("Q2").Value = 15
$ws.Range("R2").Value = 110
$ws.Range("S2").Value = 823
$ws.Range("T2").Value = 1335
$ws.Range("U2").Value = 99
$ws.Range("V2").Value = 11850
$ws.Range("W2").Value = 6
$ws.Range("X2").Value = 11885
$ws.Range("Y2").Value = 26
$ws.Range("Z2").Value = 164
$ws.Range("AA2").Value = 78
